# New crime data collected - weekly CompStat update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Volume/Number "49" -> "50" ---
$ws.Range("A8").Characters(21, 2).Text = "50"

# --- Header: report week dates ---
# Replace the later date first so the earlier edit (which changes length)
# doesn't shift the character offset of the later one.
$ws.Range("C9").Characters(47, 10).Text = "12/17/2023"
$ws.Range("C9").Characters(27, 9).Text = "12/11/2023"

# Row 14: Murder
# C14 switches from a numeric "1" to the literal text "0" (data suppressed).
$ws.Range("C14").Value = "'0"
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 13
$ws.Range("H14").Value = -46.153846153846
$ws.Range("J14").Value = 125
$ws.Range("K14").Value = -5.6
$ws.Range("L14").Value = -18.055555555555
$ws.Range("N14").Value = -75.918367346938

# Row 15: Rape
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 150
$ws.Range("F15").Value = 23
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 347
$ws.Range("J15").Value = 369
$ws.Range("K15").Value = -5.962059620596
$ws.Range("L15").Value = -3.072625698324
$ws.Range("M15").Value = 22.183098591549
$ws.Range("N15").Value = -50.286532951289

# Row 16: Robbery
$ws.Range("C16").Value = 113
$ws.Range("D16").Value = 85
$ws.Range("E16").Value = 32.941176470588
$ws.Range("F16").Value = 416
$ws.Range("G16").Value = 366
$ws.Range("H16").Value = 13.661202185792
$ws.Range("I16").Value = 4749
$ws.Range("J16").Value = 4942
$ws.Range("K16").Value = -3.905301497369
$ws.Range("L16").Value = 22.808378588052
$ws.Range("M16").Value = 6.69512469108
$ws.Range("N16").Value = -69.971546000632

# Row 17: Fel. Assault
$ws.Range("C17").Value = 137
$ws.Range("D17").Value = 123
$ws.Range("E17").Value = 11.382113821138
$ws.Range("F17").Value = 553
$ws.Range("G17").Value = 559
$ws.Range("H17").Value = -1.073345259391
$ws.Range("I17").Value = 7828
$ws.Range("J17").Value = 7127
$ws.Range("K17").Value = 9.835835554931
$ws.Range("L17").Value = 28.962108731466
$ws.Range("M17").Value = 81.245658717295
$ws.Range("N17").Value = -9.846827133479

# Row 18: Burglary
$ws.Range("C18").Value = 50
$ws.Range("D18").Value = 41
$ws.Range("E18").Value = 21.951219512195
$ws.Range("F18").Value = 205
$ws.Range("G18").Value = 172
$ws.Range("H18").Value = 19.186046511627
$ws.Range("I18").Value = 2779
$ws.Range("J18").Value = 2790
$ws.Range("K18").Value = -0.394265232974
$ws.Range("L18").Value = 29.496738117427
$ws.Range("M18").Value = -14.254859611231
$ws.Range("N18").Value = -84.79509766373

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 162
$ws.Range("D19").Value = 143
$ws.Range("E19").Value = 13.286713286713
$ws.Range("F19").Value = 667
$ws.Range("G19").Value = 607
$ws.Range("H19").Value = 9.88467874794
$ws.Range("I19").Value = 7809
$ws.Range("J19").Value = 7804
$ws.Range("K19").Value = 0.064069707842
$ws.Range("L19").Value = 16.796290756805
$ws.Range("M19").Value = 71.362738643844
$ws.Range("N19").Value = 7.119341563786

# Row 20: G.L.A.
$ws.Range("C20").Value = 93
$ws.Range("D20").Value = 75
$ws.Range("E20").Value = 24
$ws.Range("F20").Value = 305
$ws.Range("G20").Value = 328
$ws.Range("H20").Value = -7.012195121951
$ws.Range("I20").Value = 4931
$ws.Range("J20").Value = 3769
$ws.Range("K20").Value = 30.830459007694
$ws.Range("L20").Value = 66.419169760378
$ws.Range("M20").Value = 138.558297048863
$ws.Range("N20").Value = -66.772237196765

# Row 21: TOTAL
$ws.Range("C21").Value = 560
$ws.Range("D21").Value = 472
$ws.Range("E21").Value = 18.64406779661
$ws.Range("F21").Value = 2176
$ws.Range("G21").Value = 2068
$ws.Range("H21").Value = 5.22243713733
$ws.Range("I21").Value = 28561
$ws.Range("J21").Value = 26926
$ws.Range("K21").Value = 6.072197875659
$ws.Range("L21").Value = 28.456418098407
$ws.Range("M21").Value = 49.989496901586
$ws.Range("N21").Value = -56.78664911564

# Row 22: Transit
$ws.Range("C22").Value = 6
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 27
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = 42.105263157894
$ws.Range("I22").Value = 313
$ws.Range("J22").Value = 344
$ws.Range("K22").Value = -9.011627906976
$ws.Range("L22").Value = 12.589928057554
$ws.Range("M22").Value = -0.63492063492

# Row 23: Housing
$ws.Range("C23").Value = 34
$ws.Range("D23").Value = 28
$ws.Range("E23").Value = 21.428571428571
$ws.Range("F23").Value = 112
$ws.Range("G23").Value = 119
$ws.Range("H23").Value = -5.882352941176
$ws.Range("I23").Value = 1668
$ws.Range("J23").Value = 1534
$ws.Range("K23").Value = 8.735332464146
$ws.Range("L23").Value = 35.941320293398
$ws.Range("M23").Value = 60.848601735776

# Row 24: Petit Larceny
$ws.Range("C24").Value = 371
$ws.Range("D24").Value = 368
$ws.Range("E24").Value = 0.815217391304
$ws.Range("F24").Value = 1350
$ws.Range("G24").Value = 1394
$ws.Range("H24").Value = -3.156384505021
$ws.Range("I24").Value = 17278
$ws.Range("J24").Value = 17914
$ws.Range("K24").Value = -3.550295857988
$ws.Range("L24").Value = 32.836165141846
$ws.Range("M24").Value = 38.712267180475

# Row 25: Misd. Assault
$ws.Range("C25").Value = 222
$ws.Range("D25").Value = 163
$ws.Range("E25").Value = 36.196319018404
$ws.Range("F25").Value = 783
$ws.Range("G25").Value = 658
$ws.Range("H25").Value = 18.996960486322
$ws.Range("I25").Value = 10093
$ws.Range("J25").Value = 9501
$ws.Range("K25").Value = 6.23092306073
$ws.Range("L25").Value = 18.26810405437
$ws.Range("M25").Value = -5.602319491208

# Row 26: UCR Rape*
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 57.142857142857
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 29.729729729729
$ws.Range("I26").Value = 608
$ws.Range("J26").Value = 626
$ws.Range("K26").Value = -2.875399361022
$ws.Range("L26").Value = 3.050847457627

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 19
$ws.Range("D27").Value = 16
$ws.Range("E27").Value = 18.75
$ws.Range("F27").Value = 64
$ws.Range("G27").Value = 64
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 1003
$ws.Range("J27").Value = 883
$ws.Range("K27").Value = 13.590033975084
$ws.Range("L27").Value = 10.584343991179

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 38
$ws.Range("H28").Value = -60.526315789473
$ws.Range("I28").Value = 366
$ws.Range("J28").Value = 466
$ws.Range("K28").Value = -21.459227467811
$ws.Range("L28").Value = -37.328767123287
$ws.Range("M28").Value = -21.459227467811
$ws.Range("N28").Value = -73.478260869565

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = -37.5
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = 31
$ws.Range("H29").Value = -54.838709677419
$ws.Range("I29").Value = 307
$ws.Range("J29").Value = 394
$ws.Range("K29").Value = -22.081218274111
$ws.Range("L29").Value = -37.854251012145
$ws.Range("M29").Value = -22.081218274111
$ws.Range("N29").Value = -75.361155698234

# Row 30: Hate Crimes
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 27
$ws.Range("K30").Value = -35.714285714285
$ws.Range("L30").Value = -42.553191489361

